# Unhide the previously hidden detail rows on the active sheet.
# (The sheet had a block of rows hidden under each "group" header row;
#  this change makes them visible again.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsToUnhide = @(
    9, 10,
    13, 14,
    17,
    21, 22, 23, 24, 25, 26, 27, 28, 29,
    32, 33, 34,
    37, 38, 39, 40, 41, 42,
    46, 47,
    50, 51, 52,
    55, 56, 57,
    61, 62, 63,
    66, 67,
    70, 71,
    75, 76,
    79, 80,
    83, 84
)

foreach ($r in $rowsToUnhide) {
    $ws.Rows.Item($r).Hidden = $false
}
